$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the last-updated date (C1) from 2024-03-15 to 2024-03-28 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "RAF-capacity" sheet: update B24 and B25 (hydrogen combustion/combined-cycle) RAF values ---
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# --- Switch the active/visible sheet from RAF-generation to RAF-capacity, with updated selection/zoom ---
$wsCapacity.Activate() | Out-Null
$wsCapacity.Range("B25").Select() | Out-Null
$excel.ActiveWindow.Zoom = 80
